$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.590.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +6.85%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.811.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +6.11%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +4.13%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9982"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3833"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.09%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.09"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.32%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3511"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +6.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.232"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.15%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07737"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.34%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.44"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.611"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.63%  "

$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.815.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +6.64%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.195"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.68%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001124"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06700"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.06%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "86.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +6.64%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9988"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.01%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.83%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.508"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.33%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.12"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.82%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "27.564.06"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.469"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.668"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +7.60%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.08"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +15.26%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.488"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +15.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.55"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.61%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.018.75"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.78%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "136.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.19%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.397"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.77%  "

$ws.Range("E33").Value = "  -0.43%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.87"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.73%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08779"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.17%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.722"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.14%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.638"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.79%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.7086"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +15.85%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06531"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.64%  "

$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2261"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +6.92%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02413"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.16%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.989"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.293"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.60%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.97"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.47%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6605"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +13.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.9978"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.048"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.180"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.08%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.64"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07362"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.96%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.07%  "
